$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A69 precision
$ws.Range("A69").Value = 44382.76768740047

# Add new row 70 data
$ws.Range("A70").Value = 44383.76759035335
$ws.Range("B70").Value = 78651
$ws.Range("C70").Value = 66274
$ws.Range("D70").Value = 3562
$ws.Range("E70").Value = 2144
$ws.Range("F70").Value = 1526
$ws.Range("G70").Value = 20884
$ws.Range("H70").Value = 1547
$ws.Range("I70").Value = 865
$ws.Range("J70").Value = 193

# Apply same style as A69 (s="2", date format) to A70
$ws.Range("A70").NumberFormat = $ws.Range("A69").NumberFormat
